$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Api Detail")
$ws.Activate()

# --- Insert block 1: new rows 50 & 51 (after existing row 49 "Api Detail for Mto City") ---
$ws.Rows("50:51").Insert()

$ws.Range("B50").Value = "Fetch All Link Country with Mito Partner"

$ws.Range("B51").Value = "Fetch All City By Country Code"
$ws.Range("C51").Value = "Done in city contrller"
$ws.Range("C14").Copy()
$ws.Range("C51").PasteSpecial(-4122)
$ws.Range("D51").Value = "Test Done"

# --- Insert block 2: new rows 55 & 56 (after "Api Detail for Mto Wallet" header, now row 54) ---
$ws.Rows("55:56").Insert()

$ws.Range("B55").Value = "Fetch All Link Wallet with Mito Partner"

$ws.Range("B56").Value = "Fetch All Wallets By Country Code"
$ws.Range("C56").Value = "Done in country controller"
$ws.Range("C14").Copy()
$ws.Range("C56").PasteSpecial(-4122)
$ws.Range("D56").Value = "Test Done"

# --- Insert block 3: new rows 60 & 61 (after "Api Detail for Mto Bank" header, now row 59) ---
$ws.Rows("60:61").Insert()

$ws.Range("B60").Value = "Fetch All Link Bank with Mito Partner"

$ws.Range("B61").Value = "Fetch All Bank By Country Code"
$ws.Range("C61").Value = "Done in country controller"
$ws.Range("C14").Copy()
$ws.Range("C61").PasteSpecial(-4122)
$ws.Range("D61").Value = "Test Done"

$excel.CutCopyMode = 0

# --- sheet view: reflect the new selection position ---
$ws.Range("B50").Select()
